$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 8422
$ws.Range("F5").Value = 87
$ws.Range("F6").Value = 538
$ws.Range("F7").Value = 7395
$ws.Range("F8").Value = 7395
$ws.Range("F9").Value = 600
$ws.Range("F10").Value = 515
$ws.Range("F13").Value = 309
$ws.Range("F17").Value = 125
$ws.Range("F18").Value = 162
$ws.Range("F19").Value = 12269
$ws.Range("F20").Value = 107
$ws.Range("F22").Value = 2517
$ws.Range("F23").Value = 3612
$ws.Range("F26").Value = 2955
$ws.Range("F27").Value = 115
$ws.Range("F30").Value = 39
$ws.Range("F31").Value = 3367
$ws.Range("F33").Value = 344
$ws.Range("F34").Value = 1734
$ws.Range("F35").Value = 83
$ws.Range("F36").Value = 139
$ws.Range("F37").Value = 6064
$ws.Range("F38").Value = 102
$ws.Range("F39").Value = 1847
$ws.Range("F40").Value = 1261
$ws.Range("F41").Value = 36
$ws.Range("F42").Value = 919
$ws.Range("F44").Value = 172
$ws.Range("F45").Value = 9
$ws.Range("F48").Value = 1114
$ws.Range("F50").Value = 27

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 213
$ws.Range("F16").Value = 110

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 334
$ws.Range("F3").Value = 480
$ws.Range("F4").Value = 13

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 334
$ws.Range("F7").Value = 8422
$ws.Range("F8").Value = 87
$ws.Range("F10").Value = 538
$ws.Range("F11").Value = 7395
$ws.Range("F12").Value = 600
$ws.Range("F13").Value = 515
$ws.Range("F14").Value = 309
$ws.Range("F19").Value = 162
$ws.Range("F20").Value = 213
$ws.Range("F21").Value = 12269
$ws.Range("F22").Value = 107
$ws.Range("F25").Value = 2517
$ws.Range("F26").Value = 2517
$ws.Range("F27").Value = 3612
$ws.Range("F28").Value = 115
$ws.Range("F31").Value = 39
$ws.Range("F33").Value = 3367
$ws.Range("F34").Value = 344
$ws.Range("F35").Value = 1734
$ws.Range("F36").Value = 83
$ws.Range("F37").Value = 139
$ws.Range("F38").Value = 6064
$ws.Range("F40").Value = 102
$ws.Range("F41").Value = 1847
$ws.Range("F43").Value = 1261
$ws.Range("F44").Value = 36
$ws.Range("F45").Value = 919
$ws.Range("F46").Value = 172
$ws.Range("F49").Value = 1114
$ws.Range("F51").Value = 27

Write-Host "Applied all view-count updates"